$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation needs to be inserted at the top of this
# block (row 309), pushing all the existing rows for this market/category
# down by one. Use a real row insert so everything below shifts correctly.
$ws.Rows("309").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(309, 1).Value = 4
$ws.Cells.Item(309, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(309, 3).Value = "Los Lagos"
$ws.Cells.Item(309, 4).Value = 45093
$ws.Cells.Item(309, 5).Value = 10
$ws.Cells.Item(309, 6).Value = 100112044
$ws.Cells.Item(309, 7).Value = "Perejil"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 180
$ws.Cells.Item(309, 11).Value = 6000
$ws.Cells.Item(309, 12).Value = 6000
$ws.Cells.Item(309, 13).Value = 6000
$ws.Cells.Item(309, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(309, 15).Value = "Región Metropolitana"
$ws.Cells.Item(309, 16).Value = 2000
$ws.Cells.Item(309, 17).Value = 3
$ws.Cells.Item(309, 18).Value = "Hortaliza"
